$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.891.73"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "1.571.17"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("D4").Value = "'0.993"
$ws.Range("E4").Value = "  -1.35%  "

$ws.Range("D5").Value = "'211.45"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D7").Value = "'0.992"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").Value = "'23.22"
$ws.Range("E8").Value = "  +6.03%  "

$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").Value = "1.797.16"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").Value = "1.584.36"
$ws.Range("E13").Value = "  +2.33%  "

$ws.Range("D14").Value = "'3.75"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("D16").Value = "27.884.86"
$ws.Range("E16").Value = "  +3.07%  "

$ws.Range("D17").Value = "'63.31"
$ws.Range("E17").Value = "  +2.17%  "

$ws.Range("D18").Value = "'229.04"
$ws.Range("E18").Value = "  +6.48%  "

$ws.Range("D19").Value = "0.0₃0703"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").Value = "'7.44"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").Value = "'0.992"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").Value = "'4.11"
$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").Value = "'9.29"
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").Value = "'151.26"
$ws.Range("E25").Value = "  -1.79%  "

$ws.Range("D26").Value = "'15.21"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.107"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.56"
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "'0.993"
$ws.Range("E29").Value = "  -1.31%  "

$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").Value = "'3.13"
$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").Value = "1.412.16"
$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("E35").Value = "  -0.47%  "

$ws.Range("E36").Value = "  -4.32%  "

$ws.Range("E37").Value = "  -1.91%  "

$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("E40").Value = "  +2.57%  "

$ws.Range("D41").Value = "'0.805"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("E42").Value = "  -1.55%  "

$ws.Range("D43").Value = "'5.57"
$ws.Range("E43").Value = "  -3.83%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'0.973"
$ws.Range("E44").Value = "  -2.62%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.81"
$ws.Range("E45").Value = "  +4.24%  "

$ws.Range("D46").Value = "'63.72"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").Value = "1.707.51"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "'86.62"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0942"
$ws.Range("E51").Value = "  -1.65%  "
